# Apply renames / view-state updates per commit diff.
$wb = $excel.ActiveWorkbook

# Rename the two worksheets (this also updates the Print_Area defined name
# automatically, since it is a sheet-scoped name referencing Sheet1).
$wsCommands = $wb.Worksheets.Item(1)
$wsCommands.Name = "Commands"

$wsMisc = $wb.Worksheets.Item(2)
$wsMisc.Name = "Misc."

# Update the print area explicitly too, to be safe, now that the sheet has
# been renamed.
$wsCommands.PageSetup.PrintArea = '$A$1:$N$58'

# --- Misc. sheet view state ------------------------------------------------
# Move the selection from H29 to H21.
$wsMisc.Activate()
$wsMisc.Range("H21").Select() | Out-Null

# --- Commands sheet view state -------------------------------------------
# Remove the scrolled-down view (topLeftCell="A25") and move the selection
# from C31 to C71. Activate this sheet last so it remains the selected tab.
$wsCommands.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsCommands.Range("C71").Select() | Out-Null
